$d = $word.ActiveDocument

# 1. SD Kartenslot paragraph: drop "(Ausführung als 2,00 Pins)" qualifier, keep SDC0 but wrap it in parens.
$d.Content.Find.Execute(
    "SD Kartenslot (Ausführung als 2,00 Pins) SDC0",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SD Kartenslot (SDC0)", 2) | Out-Null

# 2. USB Host paragraph: drop the "Jumper 5V USB Port und " part.
$d.Content.Find.Execute(
    "USB Host (Jumper 5V USB Port und ID Pin Device / Host)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "USB Host (ID Pin Device / Host)", 2) | Out-Null

# 3. 24bit RGB Interface paragraph: reword the cabling / touch / backlight description.
$d.Content.Find.Execute(
    "24bit RGB Interface Flachbandkabel (auf kleinen Headern 0,1) Mit Touchsupport",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "24bit RGB Interface, verfügbar über ein 50 Pin FFC, mit Touchsupport (NS2009) und LED Backlight Treiber.", 2) | Out-Null

# 4. Parallel CSI paragraph: replace the whole sentence (also removes the ebay/proofErr remnants).
$d.Content.Find.Execute(
    "Parallel CSI momentan kein Softwaresupport Benutzung fertiges Modul (siehe ebay)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Parallel CSI, verfügbar über ein 24 Pin FFC Stecker. Somit besteht die Möglichkeit unterschiedliche Kameras anzuschließen.", 2) | Out-Null

# 5. MIPI CSI paragraph: reword with FFC connector detail and move the Mainline Kernel caveat inside the parenthesis.
$d.Content.Find.Execute(
    "MIPI CSI (Raspberry Pi kompatibel) Momentan kein Softwaresupport",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "MIPI CSI, verfügbar über ein 15 Pin FFC Stecker. Somit kompatibel zur Raspberry Pi Kamera (Momentan kein Softwaresupport im Mainline Kernel)", 2) | Out-Null

# 6. Audio Out / Mikrofon in paragraph: add the audio jack description.
$d.Content.Find.Execute(
    "Audio Out / Mikrofon in (Ausführung als 2,00Pins)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Audio Out / Mikrofon in über eine 4 polige 3,5mm Audio Jack (Ausführung als 2,00Pins)", 2) | Out-Null

# 7. Raspberry PI GPIO Header paragraph: prepend "26 Pin Standard ".
$d.Content.Find.Execute(
    "Raspberry PI GPIO Header (UART / I²C / SPI)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "26 Pin Standard Raspberry PI GPIO Header (UART / I²C / SPI)", 2) | Out-Null

# 8. Append three new paragraphs after "SPI Flash NOR / NAND".
$spiPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "SPI Flash NOR / NAND") {
        $spiPara = $p
    }
}

$spiPara.Range.InsertParagraphAfter()
$newPara1 = $spiPara.Next()
$newPara1.Range.Text = "Pinheader zum Anschluss einer RTC-Batterie"

$newPara1.Range.InsertParagraphAfter()
$newPara2 = $newPara1.Next()
$newPara2.Range.Text = "Einstellbare GPIO Versorgungsspannung (1,2V; 1,8V; 3,0V; 3,3V)"

$newPara2.Range.InsertParagraphAfter()
$newPara3 = $newPara2.Next()
$newPara3.Range.Text = "Verstellbare GPIO Versorgungsspannung für die GPIO Bank PE"

Write-Output "Done"
